$d = $word.ActiveDocument

function Insert-BoldHeading {
    param(
        [int]$beforeParaIndex,
        [string]$text
    )

    $target = $d.Paragraphs($beforeParaIndex)
    $target.Range.InsertParagraphBefore()

    # The freshly minted (empty) paragraph now sits at $beforeParaIndex,
    # pushing the original content down by one.
    $newPara = $d.Paragraphs($beforeParaIndex)
    $newRange = $newPara.Range
    $newRange.Text = $text

    # Bold the visible text only (exclude the paragraph mark) so the
    # paragraph's own rPr/pPr stays untouched.
    $freshPara = $d.Paragraphs($beforeParaIndex)
    $freshRange = $freshPara.Range
    $textOnly = $d.Range($freshRange.Start, $freshRange.Start + $text.Length)
    $textOnly.Font.Bold = 1

    # BoldBi (-> w:bCs) needs the full paragraph range (incl. the mark) to
    # stick reliably.
    $freshPara2 = $d.Paragraphs($beforeParaIndex)
    $freshRange2 = $freshPara2.Range
    $freshRange2.Font.BoldBi = 1
}

# 1) "Ratification of 1725 Treaty" before "Articles of Submission & Agreements..."
Insert-BoldHeading 3 "Ratification of 1725 Treaty"

# 2) "Mascarene's Promises" before "By Major Paul Mascarene..." (now shifted
#    down one slot to 13 because of the insertion above).
Insert-BoldHeading 13 "Mascarene’s Promises"
